# Updated to do list
# Insert a new task row into the "To Do- FY15 Release" sheet (row 16),
# shifting the existing rows down by one, and populate the new row with
# the new to-do item. Also move the active sheet/selection to this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do- FY15 Release")

# Insert a new blank row above the current row 16 (formatting copied from
# the row above, matching the surrounding rows' style).
$ws.Rows("16").Insert()

# Fill in the new to-do item.
$ws.Range("A16").Value = "Done"
$ws.Range("B16").Value = "Fix bug in performance ratio"
$ws.Range("C16").Value = "Janine"
$ws.Range("E16").Value = "A"

# Make this sheet the active one, with the cursor resting on the row
# below the newly-entered item (A17), matching the author's final state.
$ws.Activate()
$ws.Range("A17").Select()
